# Add two new columns, I ("I0") and J ("IF"), to the right of the existing
# H ("IP") column, matching the header's bold/bordered/centered style, then
# fill in the corresponding numeric data for rows 2-7.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone H1's formatting (bold font, thin border, centered/top alignment)
# onto the two new header cells, then set their text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data for columns I (I0) and J (IF), rows 2-7.
$iValues = @(5, 6, 7, 8, 7, 6)
$jValues = @(6, 7, 7, 8, 8, 6)

for ($idx = 0; $idx -lt $iValues.Length; $idx++) {
    $row = $idx + 2
    $ws.Cells.Item($row, 9).Value  = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
